$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Turni Fissi")

# Row 10
$rng = $ws.Range("F10:G10")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("F10").Value = "VIN"
$ws.Range("G10").Value = "MAR"

# Row 11
$rng = $ws.Range("C11:D11")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("E11").Clear()
$rng = $ws.Range("F11:G11")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("C11").Value = "CAR"
$ws.Range("D11").Value = "URG"
$ws.Range("F11").Value = "MAD"
$ws.Range("G11").Value = "CMG"

# Row 12
$rng = $ws.Range("C12:D12")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("E12").Clear()
$rng = $ws.Range("F12:G12")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("C12").Value = "VAN"
$ws.Range("D12").Value = "URG"
$ws.Range("F12").Value = "CAR"
$ws.Range("G12").Value = "VIN"

# Row 17
$rng = $ws.Range("F17:G17")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("F17").Value = "URG"
$ws.Range("G17").Value = "MAD"

# Row 18
$rng = $ws.Range("C18:D18")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("E18").Clear()
$rng = $ws.Range("F18:G18")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("C18").Value = "VIN"
$ws.Range("D18").Value = "MAR"
$ws.Range("F18").Value = "CMG"
$ws.Range("G18").Value = "CAR"

# Row 19
$rng = $ws.Range("C19:D19")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("E19").Clear()
$rng = $ws.Range("F19:G19")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("C19").Value = "EMA"
$ws.Range("D19").Value = "SAC"
$ws.Range("F19").Value = "VAN"
$ws.Range("G19").Value = "BET"

# Row 24
$rng = $ws.Range("F24:G24")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("F24").Value = "VIN"
$ws.Range("G24").Value = "URG"

# Row 25
$rng = $ws.Range("C25:D25")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("E25").Clear()
$rng = $ws.Range("F25:G25")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("C25").Value = "MAD"
$ws.Range("D25").Value = "MAR"
$ws.Range("F25").Value = "DAN"
$ws.Range("G25").Value = "SAC"

# Row 26
$rng = $ws.Range("C26:D26")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("E26").Clear()
$rng = $ws.Range("F26:G26")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("C26").Value = "VIN"
$ws.Range("D26").Value = "URG"
$ws.Range("F26").Value = "MAD"
$ws.Range("G26").Value = "MAR"

# Row 31
$rng = $ws.Range("F31:G31")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("F31").Value = "CMG"
$ws.Range("G31").Value = "CAR"

# Row 32
$rng = $ws.Range("C32:D32")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("E32").Clear()
$rng = $ws.Range("F32:G32")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("C32").Value = "EMA"
$ws.Range("D32").Value = "BET"
$ws.Range("F32").Value = "SAC"
$ws.Range("G32").Value = "VAN"

# Row 33
$rng = $ws.Range("C33:D33")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("E33").Clear()
$rng = $ws.Range("F33")
$rng.Borders.LineStyle = -4142
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$ws.Range("C33").Value = "CMG"
$ws.Range("D33").Value = "CAR"
$ws.Range("F33").Value = "BET"
$g33 = $ws.Range("G33")
$g33.Borders.LineStyle = -4142
$g33.HorizontalAlignment = -4108
$g33.VerticalAlignment = -4108
$g33.Font.Underline = 2
$g33.ClearContents()

# Selection change
$ws.Range("G33").Select()
